$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.094.55'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '2.303.46'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.83'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +1.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.80'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.30'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.64%  '
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.27'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +13.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.77'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '2.660.45'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '2.291.45'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("E18").Value = '  +3.51%  '
$ws.Range("D19").Value = '43.030.29'
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("E20").Value = '  +3.05%  '
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("E22").Value = '  +1.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.78'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.75'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("E25").Value = '  +7.04%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -0.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.45'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("E29").Value = '  +10.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.67'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.11'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.32%  '
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +0.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.66'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +7.47%  '
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.80'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +4.82%  '
$ws.Range("E38").Value = '  +0.99%  '
$ws.Range("E39").Value = '  +0.50%  '
$ws.Range("E40").Value = '  +1.70%  '
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("E42").Value = '  +0.27%  '
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = '1.983.17'
$ws.Range("E44").Value = '  +1.15%  '
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.95'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.67'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.25%  '
$ws.Range("E48").Value = '  +2.01%  '
$ws.Range("D49").Value = '2.529.34'
$ws.Range("E49").Value = '  +1.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.33'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.59'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.58%  '
